$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$authors = @(
    '[Fei%Zhou%NULL%0,                  Ting%Yu%NULL%0,                  Ronghui%Du%NULL%0,                  Guohui%Fan%NULL%0,                  Ying%Liu%NULL%0,                  Zhibo%Liu%NULL%0,                  Jie%Xiang%NULL%0,                  Yeming%Wang%NULL%0,                  Bin%Song%NULL%0,                  Xiaoying%Gu%NULL%0,                  Lulu%Guan%NULL%0,                  Yuan%Wei%NULL%0,                  Hui%Li%NULL%0,                  Xudong%Wu%NULL%0,                  Jiuyang%Xu%NULL%0,                  Shengjin%Tu%NULL%0,                  Yi%Zhang%NULL%0,                  Hua%Chen%NULL%0,                  Bin%Cao%NULL%0]',
    '[Carb\u00f3-Cisnero%Yacquelin%coreGivesNoEmail%0,                 Fern\u00e1ndez-Gonz\u00e1lez%Paula%coreGivesNoEmail%0,                 Hierrezuelo-Rojas%Naifi%coreGivesNoEmail%0,                 Subert-Salas%Lizandra%coreGivesNoEmail%0]',
    '[Kyung Soo%Hong%NULL%0,                  Kwan Ho%Lee%NULL%2,                  Kwan Ho%Lee%NULL%0,                  Jin Hong%Chung%NULL%2,                  Jin Hong%Chung%NULL%0,                  Kyeong-Cheol%Shin%NULL%2,                  Kyeong-Cheol%Shin%NULL%0,                  Eun Young%Choi%NULL%2,                  Eun Young%Choi%NULL%0,                  Hyun Jung%Jin%NULL%2,                  Hyun Jung%Jin%NULL%0,                  Jong Geol%Jang%NULL%2,                  Jong Geol%Jang%NULL%0,                  Wonhwa%Lee%NULL%2,                  Wonhwa%Lee%NULL%0,                  June Hong%Ahn%NULL%0,                  June Hong%Ahn%NULL%0]',
    '[Matt%Arentz%xref no email%0,           Eric%Yim%xref no email%0,           Lindy%Klaff%xref no email%0,           Sharukh%Lokhandwala%xref no email%0,           Francis X.%Riedo%xref no email%0,           Maria%Chong%xref no email%0,           Melissa%Lee%xref no email%0]',
    '[Aggarwal%Gaurav%coreGivesNoEmail%0,                 Aggarwal%Saurabh%coreGivesNoEmail%0,                 Garcia-Telles%Nelson%coreGivesNoEmail%0,                 Henry%Brandon Michael%coreGivesNoEmail%0,                 Lavie%Carl%coreGivesNoEmail%0,                 Lippi%Giuseppe%coreGivesNoEmail%0]',
    '[Yuan%Yu%NULL%0,                  Dan%Xu%NULL%0,                  Shouzhi%Fu%NULL%0,                  Jun%Zhang%NULL%0,                  Xiaobo%Yang%NULL%0,                  Liang%Xu%NULL%0,                  Jiqian%Xu%NULL%0,                  Yongran%Wu%NULL%0,                  Chaolin%Huang%NULL%0,                  Yaqi%Ouyang%NULL%0,                  Luyu%Yang%NULL%0,                  Minghao%Fang%NULL%0,                  Hongwen%Xiao%NULL%0,                  Jing%Ma%NULL%0,                  Wei%Zhu%NULL%0,                  Song%Hu%NULL%0,                  Quan%Hu%NULL%0,                  Daoyin%Ding%NULL%0,                  Ming%Hu%NULL%0,                  Guochao%Zhu%NULL%0,                  Weijiang%Xu%NULL%0,                  Jun%Guo%NULL%0,                  Jinglong%Xu%NULL%0,                  Haitao%Yuan%NULL%0,                  Bin%Zhang%NULL%0,                  Zhui%Yu%yuzhui@whu.edu.cn%0,                  Dechang%Chen%icudechangchen@163.com%0,                  Shiying%Yuan%yuan_shiying@163.com%0,                  You%Shang%you_shanghust@163.com%0]',
    '[Fan%Yang%NULL%0,                  Shaobo%Shi%NULL%2,                  Shaobo%Shi%NULL%0,                  Jiling%Zhu%NULL%1,                  Jinzhi%Shi%NULL%1,                  Kai%Dai%NULL%1,                  Xiaobei%Chen%121278395@qq.com%1]',
    '[Jia-Fu%Wei%NULL%0,                  Fang-Yang%Huang%NULL%0,                  Tian-Yuan%Xiong%NULL%0,                  Tian-Yuan%Xiong%NULL%0,                  Qi%Liu%NULL%0,                  Hong%Chen%NULL%0,                  Hui%Wang%NULL%0,                  He%Huang%NULL%0,                  Yi-Chun%Luo%NULL%0,                  Xuan%Zhou%NULL%0,                  Zhi-Yue%Liu%NULL%0,                  Yong%Peng%NULL%0,                  Yuan-Ning%Xu%NULL%0,                  Bo%Wang%NULL%0,                  Ying-Ying%Yang%NULL%0,                  Zong-An%Liang%NULL%0,                  Xue-Zhong%Lei%NULL%0,                  Yang%Ge%NULL%0,                  Ming%Yang%NULL%0,                  Ling%Zhang%NULL%0,                  Ming-Quan%Zeng%NULL%0,                  He%Yu%NULL%0,                  Kai%Liu%NULL%0,                  Yu-Heng%Jia%NULL%0,                  Bernard D%Prendergast%NULL%0,                  Wei-Min%Li%NULL%0,                  Wei-Min%Li%NULL%0,                  Mao%Chen%NULL%0]',
    '[Xiaochen%Li%NULL%0,                  Shuyun%Xu%NULL%0,                  Muqing%Yu%NULL%0,                  Ke%Wang%NULL%0,                  Yu%Tao%NULL%0,                  Ying%Zhou%NULL%0,                  Jing%Shi%NULL%0,                  Min%Zhou%NULL%0,                  Bo%Wu%NULL%0,                  Zhenyu%Yang%NULL%0,                  Cong%Zhang%NULL%0,                  Junqing%Yue%NULL%0,                  Zhiguo%Zhang%NULL%0,                  Harald%Renz%NULL%0,                  Xiansheng%Liu%NULL%0,                  Jungang%Xie%NULL%0,                  Min%Xie%NULL%0,                  Jianping%Zhao%NULL%0]',
    '[Anuradha%Lala%NULL%1,                  Kipp W.%Johnson%NULL%1,                  James L.%Januzzi%NULL%1,                  Adam J.%Russak%NULL%1,                  Ishan%Paranjpe%NULL%0,                  Felix%Richter%NULL%0,                  Shan%Zhao%NULL%0,                  Sulaiman%Somani%NULL%0,                  Tielman%Van Vleck%NULL%1,                  Akhil%Vaid%NULL%0,                  Fayzan%Chaudhry%NULL%1,                  Jessica K.%De Freitas%NULL%1,                  Zahi A.%Fayad%NULL%1,                  Sean P.%Pinney%NULL%1,                  Matthew%Levin%NULL%1,                  Alexander%Charney%NULL%1,                  Emilia%Bagiella%NULL%0,                  Jagat%Narula%NULL%0,                  Benjamin S.%Glicksberg%NULL%1,                  Girish%Nadkarni%NULL%1,                  Donna M.%Mancini%NULL%1,                  Valentin%Fuster%NULL%0,                  NULL%NULL%NULL%0]',
    '[Chaolin%Huang%NULL%0,                  Yeming%Wang%NULL%0,                  Xingwang%Li%NULL%0,                  Lili%Ren%NULL%0,                  Jianping%Zhao%NULL%0,                  Yi%Hu%NULL%0,                  Li%Zhang%NULL%0,                  Guohui%Fan%NULL%0,                  Jiuyang%Xu%NULL%0,                  Xiaoying%Gu%NULL%0,                  Zhenshun%Cheng%NULL%0,                  Ting%Yu%NULL%0,                  Jiaan%Xia%NULL%0,                  Yuan%Wei%NULL%0,                  Wenjuan%Wu%NULL%0,                  Xuelei%Xie%NULL%0,                  Wen%Yin%NULL%0,                  Hui%Li%NULL%0,                  Min%Liu%NULL%0,                  Yan%Xiao%NULL%0,                  Hong%Gao%NULL%0,                  Li%Guo%NULL%0,                  Jungang%Xie%NULL%0,                  Guangfa%Wang%NULL%0,                  Rongmeng%Jiang%NULL%0,                  Zhancheng%Gao%NULL%0,                  Qi%Jin%NULL%0,                  Jianwei%Wang%wangjw28@163.com%0,                  Bin%Cao%caobin_ben@163.com%0]',
    '[Huan%Han%NULL%0,                  Linlin%Xie%NULL%1,                  Rui%Liu%NULL%1,                  Jie%Yang%NULL%0,                  Fang%Liu%NULL%2,                  Kailang%Wu%NULL%2,                  Lang%Chen%NULL%1,                  Wei%Hou%NULL%0,                  Yong%Feng%yongfeng@whu.edu.cn%0,                  Chengliang%Zhu%xinchengzhu@163.com%2,                  Chengliang%Zhu%xinchengzhu@163.com%0]',
    '[Qing%Deng%NULL%0,                  Bo%Hu%NULL%1,                  Yao%Zhang%NULL%0,                  Hao%Wang%NULL%0,                  Xiaoyang%Zhou%NULL%1,                  Wei%Hu%NULL%1,                  Yuting%Cheng%NULL%1,                  Jie%Yan%NULL%0,                  Haiqin%Ping%NULL%1,                  Qing%Zhou%qingzhou.wh.edu@hotmail.com%1]',
    '[Rongrong%Yang%NULL%1,                  Xien%Gui%NULL%1,                  Yongxi%Zhang%NULL%0,                  Yong%Xiong%NULL%0]',
    '[Xiaobo%Yang%NULL%0,                  Yuan%Yu%NULL%0,                  Jiqian%Xu%NULL%0,                  Huaqing%Shu%NULL%0,                  Jia''an%Xia%NULL%0,                  Hong%Liu%NULL%0,                  Yongran%Wu%NULL%0,                  Lu%Zhang%NULL%0,                  Zhui%Yu%NULL%0,                  Minghao%Fang%NULL%0,                  Ting%Yu%NULL%0,                  Yaxin%Wang%NULL%0,                  Shangwen%Pan%NULL%0,                  Xiaojing%Zou%NULL%0,                  Shiying%Yuan%NULL%0,                  You%Shang%NULL%0]',
    '[Aguilar%Medina Jos\u00e9 M.%coreGivesNoEmail%0,                 Claro%Vald\u00e9s Ram\u00f3n%coreGivesNoEmail%1,                 Obreg\u00f3n%Santos Angel G.%coreGivesNoEmail%1,                 Proh\u00edas%Mart\u00ednez Juan%coreGivesNoEmail%1,                 Rodr\u00edguez%Blanco Suilbert%coreGivesNoEmail%1]',
    '[Shao-Fang%Nie%NULL%0,                  Miao%Yu%NULL%0,                  Tian%Xie%NULL%0,                  Fen%Yang%NULL%0,                  Hong-Bo%Wang%NULL%0,                  Zhao-Hui%Wang%NULL%0,                  Ming%Li%NULL%0,                  Xing-Li%Gao%NULL%0,                  Bing-Jie%Lv%NULL%0,                  Shi-Jia%Wang%NULL%0,                  Xiao-Bo%Zhang%NULL%0,                  Shao-Lin%He%NULL%0,                  Zhi-Hua%Qiu%NULL%0,                  Yu-Hua%Liao%NULL%0,                  Zi-Hua%Zhou%NULL%0,                  Xiang%Cheng%NULL%0]',
    '[Tao%Guo%NULL%0,                  Yongzhen%Fan%NULL%0,                  Ming%Chen%NULL%0,                  Xiaoyan%Wu%NULL%0,                  Lin%Zhang%NULL%0,                  Tao%He%NULL%0,                  Hairong%Wang%NULL%0,                  Jing%Wan%NULL%0,                  Xinghuan%Wang%NULL%0,                  Zhibing%Lu%NULL%0]',
    '[Yan%Deng%NULL%0,                  Wei%Liu%NULL%0,                  Kui%Liu%NULL%0,                  Yuan-Yuan%Fang%NULL%0,                  Jin%Shang%NULL%1,                  Ling%Zhou%NULL%0,                  Ke%Wang%NULL%0,                  Fan%Leng%NULL%1,                  Shuang%Wei%NULL%0,                  Lei%Chen%NULL%1,                  Hui-Guo%Liu%NULL%0,                  Pei-Fang%Wei%NULL%0,                  Pei-Fang%Wei%NULL%0]',
    '[Tao%Chen%NULL%0,                  Di%Wu%NULL%0,                  Huilong%Chen%NULL%0,                  Weiming%Yan%NULL%0,                  Danlei%Yang%NULL%0,                  Guang%Chen%NULL%0,                  Ke%Ma%NULL%0,                  Dong%Xu%NULL%0,                  Haijing%Yu%NULL%0,                  Hongwu%Wang%NULL%0,                  Tao%Wang%NULL%0,                  Wei%Guo%NULL%0,                  Jia%Chen%NULL%0,                  Chen%Ding%NULL%0,                  Xiaoping%Zhang%NULL%0,                  Jiaquan%Huang%NULL%0,                  Meifang%Han%NULL%0,                  Shusheng%Li%NULL%0,                  Xiaoping%Luo%NULL%0,                  Jianping%Zhao%NULL%0,                  Qin%Ning%NULL%0]'
)

for ($i = 0; $i -lt $authors.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $authors[$i]
}

Write-Output "Done"